# Update CDA Logical model for ST.r2b
# - Bump Version / Date metadata
# - Add a new "Jurisdiction" property row to the Metadata sheet
# - Add the II-1 constraint text to SubjectPerson.typeId (Elements sheet)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Metadata sheet
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.0.0-sd-202406-matchbox-patch -> 2.0.1-sd-202510-matchbox-patch
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Date: 2024-06-19T17:47:42+02:00 -> 2025-10-29T22:15:57+01:00
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" row right after "Contact" (row 10) and before
# "Description" (row 11), pushing every following row down by one.
$meta.Rows.Item(11).Insert()

# Re-apply the same formatting used by the rest of the property rows so the
# new row matches the existing style (instead of the default style Excel
# assigns to a freshly inserted row).
$meta.Range("A12:B12").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# ---------------------------------------------------------------------------
# Elements sheet
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# SubjectPerson.typeId (row 5) gains a new Constraint(s) entry (column AJ).
$elements.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}`n"
